$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 36
$ws.Range("E4").Value = 14
$ws.Range("E9").Value = 18
$ws.Range("E15").Value = 130
$ws.Range("E17").Value = 84
$ws.Range("F17").Value = 37
$ws.Range("H17").Value = 37
$ws.Range("E18").Value = 82
$ws.Range("E19").Value = 34
$ws.Range("E26").Value = 19
$ws.Range("E36").Value = 68
$ws.Range("E37").Value = 33
$ws.Range("E40").Value = 11
$ws.Range("F40").Value = 6
$ws.Range("H40").Value = 6
$ws.Range("E49").Value = 53
$ws.Range("E60").Value = 12
$ws.Range("E62").Value = 29
$ws.Range("F62").Value = 5
$ws.Range("H62").Value = 5
$ws.Range("E66").Value = 27
$ws.Range("E67").Value = 30
$ws.Range("E72").Value = 27
$ws.Range("E73").Value = 19
$ws.Range("E77").Value = 37
$ws.Range("F77").Value = 14
$ws.Range("H77").Value = 14
$ws.Range("E81").Value = 9
$ws.Range("E88").Value = 14
